$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2 ("description"): update header text of A1 ---
$ws2.Range("A1").Value = "เลขที่  (null)  ตัวเลขเท่านั้น"

# --- Sheet2: add new columns F & G with a merged, centered header ---
$ws2.Range("F1:G1").HorizontalAlignment = -4108  # xlCenter
$ws2.Range("F1:G1").Merge()

$ws2.Columns.Item(6).ColumnWidth = 25.6
$ws2.Columns.Item(7).ColumnWidth = 25.6

# --- Sheet2: new descriptive cells under the new columns ---
$ws2.Range("F2").Value = "เป็นค่าว่างได้"
$ws2.Range("G2").Value = "เพิ่มข้อมูล ต้องไม่เป็นค่าว่าง"

# Reuse the same look as the existing "Neutral"/"Bad" styled cells (A2/B2)
$ws2.Range("A2").Copy()
$ws2.Range("F2").PasteSpecial(-4122)  # xlPasteFormats

$ws2.Range("B2").Copy()
$ws2.Range("G2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Sheet1: it is no longer the active tab; selection moves to a full-column pick ---
$ws1.Range("F1:F1048576").Select()

# --- Sheet2 becomes the active tab, with a new active cell ---
$ws2.Activate()
$ws2.Range("B8").Select()
